# chgt xxx colonne export + selection ihm
#
# - Rename the "Chapitre" export header in column B (row 1) to
#   "Chapitrezzzzzzzzzzz".
# - Move the active IHM selection to cell B1 (was AC2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Exigences")

$ws.Range("B1").Value = "Chapitrezzzzzzzzzzz"

[void]$ws.Range("B1").Select()
